# PFRS submission template: re-add multi line char fix (convCRLF) after
# conflict, plus swap the "Other Owned Parcels" show/hide block from a
# showBegin/showEnd pair (keyed off hasOtherParcelsInCommunity) to a
# hideBegin/hideEnd pair keyed off otherParcelsDescription, and tidy up a
# couple of stray "s" runs / table column widths.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Swap the section-2 showBegin marker for a hideBegin marker keyed off
#    otherParcelsDescription instead of hasOtherParcelsInCommunity.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "{d.hasOtherParcelsInCommunity:ifEQ(true):showBegin}", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "{d.otherParcelsDescription:ifEM():hideBegin}", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Add convCRLF to the otherParcelsDescription show(.noData) merge tag.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "{d.otherParcelsDescription:ifEM():show(.noData)}", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "{d.otherParcelsDescription:convCRLF:ifEM():show(.noData)}", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Insert a new paragraph right after the otherParcelsDescription cell's
#    text with the matching hideEnd marker (orange, non-bold run).
# ---------------------------------------------------------------------
$otherParcelsTable = $d.Tables.Item(8)
$descCell = $otherParcelsTable.Cell(1, 2)
$cellRange = $descCell.Range
$insertStart = $cellRange.End
$cellRange.Collapse(0)
$cellRange.InsertAfter([char]13 + "{d.otherParcelsDescription:ifEM():hideEnd}")

$newRunRange = $d.Range($insertStart, $d.Content.End)
$newRunRange.Find.Execute(
    "{d.otherParcelsDescription:ifEM():hideEnd}", $false, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$newRunRange.Font.Bold = $false
$newRunRange.Font.BoldBi = $false
$newRunRange.Font.Color = 23736

# ---------------------------------------------------------------------
# 4. Remove the old showEnd run (and its trailing space) that used to
#    close the hasOtherParcelsInCommunity conditional block.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "{d.hasOtherParcelsInCommunity:ifEQ(true):showEnd} ", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. Add convCRLF to the remaining plain ifEM():show(.noData) tags that
#    are re-introduced by this fix.
# ---------------------------------------------------------------------
$simpleConvCRLFFields = @(
    "parcelsAgricultureDescription",
    "parcelsAgricultureImprovementDescription",
    "parcelsNonAgricultureUseDescription",
    "purpose",
    "soilFillTypeToPlace",
    "soilTypeRemoved",
    "soilAlternativeMeasures"
)
foreach ($field in $simpleConvCRLFFields) {
    $old = "{d.$field" + ":ifEM():show(.noData)}"
    $new = "{d.$field" + ":convCRLF:ifEM():show(.noData)}"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 6. soilReduceNegativeImpacts gets convCRLF too, and also absorbs (and
#    drops) the stray trailing "s" run that followed it.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "{d.soilReduceNegativeImpacts:ifEM():show(.noData)}s", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "{d.soilReduceNegativeImpacts:convCRLF:ifEM():show(.noData)}", 2) | Out-Null

# ---------------------------------------------------------------------
# 7. Drop the other stray trailing "s" run after the "What steps..."
#    question text in the same table.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "What steps will be taken to reduce impacts to surrounding agricultural land?s",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "What steps will be taken to reduce impacts to surrounding agricultural land?", 2) | Out-Null

# ---------------------------------------------------------------------
# 8. Nudge the two "Soil/Fill to be Removed/Placed" tables' outer column
#    widths by 1 twip each (3638->3637, 3624->3625); middle column
#    (3629) stays put.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Columns.Count -eq 3) {
        $firstWidth = $t.Columns.Item(1).Width
        if ([Math]::Abs($firstWidth - 181.9) -lt 0.05) {
            $t.Columns.Item(1).Width = 181.85
            $t.Columns.Item(3).Width = 181.25
        }
    }
}
